$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 51: Stellar -> Arweave (name/link change) ---
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"

# --- Price (column D) updates ---
# Some new price strings parse as plain numbers under Excel's default
# "General" cell format (e.g. "183.24"), which would silently turn the
# cell into a numeric cell and corrupt the exact text representation the
# source data relies on (trailing zeros, etc). Force those cells to Text
# format first so the assigned string is stored verbatim, matching the
# original inline-string cells.
$ws.Range("D2").Value = "66.907.71"
$ws.Range("D3").Value = "3.342.35"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.24"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.405"
$ws.Range("D12").Value = "3.921.73"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.23"
$ws.Range("D15").Value = "66.944.66"
$ws.Range("D17").Value = "3.337.60"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "435.91"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.71"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.68"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.79"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.11"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.90"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.35"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.03"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.49"
$ws.Range("D39").Value = "2.840.84"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.795"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.47"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.26"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0679"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.30"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.64"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "324.64"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.995"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "30.99"

# --- Volume(1h) (column E) updates ---
# These percentage strings keep leading/trailing spaces, so Excel never
# mistakes them for numeric values; plain assignment is safe.
$ws.Range("E2").Value = "  -3.74%  "
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("E6").Value = "  -4.44%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("E9").Value = "  -3.62%  "
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("E11").Value = "  -3.93%  "
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("E14").Value = "  -4.93%  "
$ws.Range("E15").Value = "  -3.72%  "
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -5.02%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("E31").Value = "  -4.76%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  -2.89%  "
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("E38").Value = "  -4.14%  "
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("E41").Value = "  -2.98%  "
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("E45").Value = "  -3.32%  "
$ws.Range("E46").Value = "  -6.17%  "
$ws.Range("E47").Value = "  -6.01%  "
$ws.Range("E48").Value = "  -4.15%  "
$ws.Range("E49").Value = "  -4.44%  "
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("E51").Value = "  -5.65%  "
